# Weekly update: insert two new rows of "Tomate" price data (week of 2021-09-09,
# serial 44448) at the top of the data block for Femacal de La Calera, pushing
# the existing rows 612-624 down to 614-626.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 612 (existing rows 612:624 shift down to 614:626).
$ws.Rows("612:613").Insert()

# New row 612: Segunda, Región de Arica y Parinacota
$ws.Range("A612").Value = 3
$ws.Range("B612").Value = "Femacal de La Calera"
$ws.Range("C612").Value = "Coquimbo"
$ws.Range("D612").Value = 44448
$ws.Range("E612").Value = 5
$ws.Range("F612").Value = 100112020
$ws.Range("G612").Value = "Tomate"
$ws.Range("H612").Value = "Larga vida"
$ws.Range("I612").Value = "Segunda"
$ws.Range("J612").Value = 410
$ws.Range("K612").Value = 18500
$ws.Range("L612").Value = 19000
$ws.Range("M612").Value = 18780
$ws.Range("N612").Value = "$/bandeja 18 kilos"
$ws.Range("O612").Value = "Región de Arica y Parinacota"
$ws.Range("P612").Value = 1043
$ws.Range("Q612").Value = 18
$ws.Range("R612").Value = "Hortaliza"

# New row 613: Tercera, Región de Arica y Parinacota
$ws.Range("A613").Value = 3
$ws.Range("B613").Value = "Femacal de La Calera"
$ws.Range("C613").Value = "Coquimbo"
$ws.Range("D613").Value = 44448
$ws.Range("E613").Value = 5
$ws.Range("F613").Value = 100112020
$ws.Range("G613").Value = "Tomate"
$ws.Range("H613").Value = "Larga vida"
$ws.Range("I613").Value = "Tercera"
$ws.Range("J613").Value = 270
$ws.Range("K613").Value = 14000
$ws.Range("L613").Value = 17000
$ws.Range("M613").Value = 15556
$ws.Range("N613").Value = "$/bandeja 18 kilos"
$ws.Range("O613").Value = "Región de Arica y Parinacota"
$ws.Range("P613").Value = 864
$ws.Range("Q613").Value = 18
$ws.Range("R613").Value = "Hortaliza"
